$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "85÷5="
$t.Cell(1, 2).Range.Text = "62÷8="
$t.Cell(1, 3).Range.Text = "36÷9="
$t.Cell(1, 4).Range.Text = "82÷5="
$t.Cell(1, 5).Range.Text = "76÷4="
$t.Cell(5, 1).Range.Text = "51÷8="
$t.Cell(5, 2).Range.Text = "47÷6="
$t.Cell(5, 3).Range.Text = "89÷2="
$t.Cell(5, 4).Range.Text = "80÷4="
$t.Cell(5, 5).Range.Text = "17÷8="
$t.Cell(9, 1).Range.Text = "33÷7="
$t.Cell(9, 2).Range.Text = "92÷6="
$t.Cell(9, 3).Range.Text = "35÷6="
$t.Cell(9, 4).Range.Text = "28÷7="
$t.Cell(9, 5).Range.Text = "85÷2="
$t.Cell(13, 1).Range.Text = "96÷2="
$t.Cell(13, 2).Range.Text = "68÷2="
$t.Cell(13, 3).Range.Text = "17÷7="
$t.Cell(13, 4).Range.Text = "54÷6="
$t.Cell(13, 5).Range.Text = "97÷7="
$t.Cell(17, 1).Range.Text = "84÷8="
$t.Cell(17, 2).Range.Text = "64÷3="
$t.Cell(17, 3).Range.Text = "36÷3="
$t.Cell(17, 4).Range.Text = "58÷2="
$t.Cell(17, 5).Range.Text = "13÷4="
